## Notification-Demo-Plan.docx edit
## "Mit Bild" list item: append " implementiert" (highlighted green)
## after the existing "?" placeholder, matching the style already used
## for "Standard implementiert" / "Progress-Bar implementiert".

$d = $word.ActiveDocument

# Locate the "Mit Bild" list paragraph (currently reads "Mit Bild<TAB>?").
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Mit Bild*?*") {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    throw "Could not find the 'Mit Bild' paragraph"
}

# Append a plain space run right before the paragraph mark.
$target.Range.InsertAfter(" ")

# Append the "implementiert" run right before the paragraph mark.
# (Re-fetch the range so it reflects the text just inserted.)
$target.Range.InsertAfter("implementiert")

# Highlight only the newly added "implementiert" word in bright green,
# i.e. OOXML <w:highlight w:val="green"/>.
$wordRange = $target.Range
$found = $wordRange.Find.Execute("implementiert", $true, $false, $false, $false, $false, `
                                  $true, 1, $false, "", 0)
if ($found) {
    $wordRange.Font.HighlightColorIndex = 4  # wdBrightGreen -> w:highlight val="green"
}
